$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Sheet "user" -- add a space before the trailing "?" in every question
# ---------------------------------------------------------------------------
$wsUser = $wb.Worksheets.Item("user")

$wsUser.Range("B1").Value = "Quante fattorie didattiche ci sono in provincia di [PROVINCIA](Salerno) ?"
$wsUser.Range("B2").Value = "Ci sono fattorie didattiche in provincia di [PROVINCIA](Salerno) ?"
$wsUser.Range("B3").Value = "Ci sono fattorie didattiche nella provincia di [PROVINCIA](Napoli) ?"
$wsUser.Range("B4").Value = "Ci sono fattorie didattiche ad [SEDE](Eboli) ?"
$wsUser.Range("B5").Value = "Quante fattorie didattiche ci sono a [SEDE](Salerno) ?"
$wsUser.Range("B6").Value = "Quante fattorie didattiche ci sono in provincia di [PROVINCIA](Salerno) con [categoria](allevamento) ?"
$wsUser.Range("B7").Value = "Ci sono fattorie didattiche in provincia di [PROVINCIA](Salerno) con produzione di [CATEGORIA](formaggi) ?"
$wsUser.Range("B8").Value = "Ci sono fattorie didattiche nella provincia di [PROVINCIA](Napoli) che producono [CATEGORIA](formaggi) ?"
$wsUser.Range("B9").Value = "Ci sono fattorie didattiche nella provincia di [PROVINCIA](Napoli) con coltivazione di [CATEGORIA](frutta) ?"

# ---------------------------------------------------------------------------
# 2) Sheet "entities_slots" -- append new city entries below the existing list
# ---------------------------------------------------------------------------
$wsSlots = $wb.Worksheets.Item("entities_slots")

$wsSlots.Range("A23").Value = "CALVI RISORTA"
$wsSlots.Range("A24").Value = "PIANO DI SORRENTO"
$wsSlots.Range("A25").Value = "CAPUA"
$wsSlots.Range("A26").Value = "MARIGLIANO"
$wsSlots.Range("A27").Value = "FRANCOLISE"
$wsSlots.Range("A28").Value = "SESSA AURUNCA"
$wsSlots.Range("A29").Value = "ASCEA MARINA"
$wsSlots.Range("A30").Value = "PIGNATARO MAGGIORE"
$wsSlots.Range("A31").Value = "CASTELNUOVO CILENTO"
$wsSlots.Range("A32").Value = "CASTEL CAMPAGNANO"

# A handful of rows in this last batch (24, 30-33) were entered in a separate
# pass and picked up a distinct (but visually identical) cell style -- carry
# that style across onto those cells, including the trailing blank row 33.
$pastedStyle = $wsSlots.Range("A23").Style
$wsSlots.Range("A24").Style = $pastedStyle
$wsSlots.Range("A30").Style = $pastedStyle
$wsSlots.Range("A31").Style = $pastedStyle
$wsSlots.Range("A32").Style = $pastedStyle
$wsSlots.Range("A33").Style = $pastedStyle

# ---------------------------------------------------------------------------
# 3) View / selection state left by the edit session
# ---------------------------------------------------------------------------
$wsSlots.Activate()
$wsSlots.Range("A31").Select()

$excel.ActiveWindow.ScrollRow = 25

$wsUser.Select()
$wsUser.Range("B10").Select()

$wsSlots.Activate()
